# Replace the digit labels ("1".."300") used for the English/French number
# columns with their spelled-out word forms, leaving the Fulfulde column (C)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$english = @(
    'one', 'two ', 'three', 'four', 'five', 'six', 'seven', 'eight', 'nine', 'ten',
    'elleven', 'twelve', 'thirteen', 'fourteen', 'fifteen', 'sixteen', 'seventeen',
    'eighteen', 'nineteen', 'twenty', 'twenty-one', 'twenty-two', 'twenty-three',
    'thirty', 'forty', 'fifty', 'sixty', 'seventy', 'eighty', 'ninety',
    'one-hundred', 'two-hundred', 'three-hundred'
)

$french = @(
    'un', 'deux', 'trois', 'quatre', 'cinq', 'six', 'sept', 'huit', 'neuf', 'dix',
    'onze', 'douze', 'treize', 'quatorze', 'quinze', 'seize', 'dix-sept',
    'dix-huit', 'dix-neuf', 'vingt', 'vingt et un', 'vingt-deux', 'vingt-trois',
    'trente', 'quarante', 'cinquante', 'soixante', 'soixante-dix', 'quatre-vingts',
    'quatre-vingt-dix', 'cent', 'deux cents', 'trois cents'
)

$startRow = 374
for ($i = 0; $i -lt $english.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $english[$i]
}
for ($i = 0; $i -lt $french.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $french[$i]
}

$wb.Save()
